$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")
$ws.Activate()

# Updated Masterdata as per 2nd may Data Refresh
# reg_center_user id 10002 -> 10003 (row 3)
$ws.Range("A3").Value = 10003
# reg_center_user id 10005 -> 10003 (row 25)
$ws.Range("A25").Value = 10003

# Restore the saved view/selection state of the sheet (scrolled/selected cell
# when the workbook was last saved)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
